# Doing Updates for Financials
# Update figures on the "ESNC" sheet (Income Statement section):
#   Row 21 - Earnings Before Interest And Taxes
#   Row 48 - Income Tax Expense
#   Row 49 - Minority Interest

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: Earnings Before Interest And Taxes
$ws.Range("D21").Value = -13000
$ws.Range("E21").Value = -3800
$ws.Range("F21").Value = -17400
$ws.Range("G21").Value = -12100
$ws.Range("H21").Value = -8300
$ws.Range("I21").Value = -10900
$ws.Range("J21").Value = -11000

# Row 48: Income Tax Expense
$ws.Range("D48").Value = 800
$ws.Range("E48").Value = 3400
$ws.Range("G48").Value = 8300

# Row 49: Minority Interest
$ws.Range("D49").Value = 1900
$ws.Range("E49").Value = 1000
$ws.Range("G49").Value = 900
